$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the timestamp precision on the existing last row (row 22)
$ws.Range("A22").Value = 45874.87520287037

# Append the new reading as row 23
$ws.Range("A23").Value = 45874.91690624816
$ws.Range("B23").Value = 2025
$ws.Range("C23").Value = 19
$ws.Range("D23").Value = 13.72
$ws.Range("E23").Value = 90.22
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0.15
$ws.Range("H23").Value = "E"
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = "22:00:20"

# Match the style used for timestamps in column A
$ws.Range("A23").NumberFormat = "YYYY-MM-DD HH:MM:SS"
